# "New Submission Synced" - append the newly received form submission as
# row 11 on the "JSS 3D" results sheet (dimension grows from A1:D10 to
# A1:D11).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$row = 11

# Timestamp and Full Name are plain text - same as every other row in
# this sheet - so a normal .Value assignment is fine for them.
$ws.Cells.Item($row, 1).Value = "2026-02-10 22:35:31"
$ws.Cells.Item($row, 2).Value = "Audu Ali Abubakar "

# Admission No ("7") is stored as text in this sheet (see e.g. C3="38",
# C4="24", C7="1"), but it looks like a plain number, so Excel would
# normally auto-convert it. Force text entry by switching the cell to
# the Text format before typing the value, then drop the cell style
# back to Normal so the stored value keeps its text type without
# leaving a stray number-format style behind.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "7"
$ws.Cells.Item($row, 3).Style = "Normal"

# AI Score is numeric in every row.
$ws.Cells.Item($row, 4).Value = 7
